# Automatische test-sync: 2025-08-28 21:05:50
# Adds a new "Opvolging retour" log entry (row 20) to the Logs sheet and
# bumps the "Retour / Terugbetaling" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# --- New row 20 on the Logs sheet ---------------------------------------
$logs.Range("A20").Value = "Opvolging retour"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("D20").Value = "Retour / Terugbetaling"
$logs.Range("F20").Value = "2025-08-28 21:05:09"
$logs.Range("G20").Value = "Nee"
$logs.Range("H20").Value = "Ja"
$logs.Range("I20").Value = "Nee"
$logs.Range("J20").Value = "Nee"

# --- Extend the conditional-formatting ranges to include the new row ---
$logs.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))
$logs.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))
$logs.Range("H2:H19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H20"))
$logs.Range("I2:I19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I20"))
$logs.Range("J2:J19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J20"))

# --- Bump the Dashboard tally for "Retour / Terugbetaling" -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 18
